# Applies two changes to the deck:
#  1. Slide 16's table switches from table style {BDACEECD-2F8F-42FE-80BF-FA8F8758E636}
#     to table style {4357684D-0C54-446E-AAA6-5CBF29CCA078}.
#  2. The presentation theme's colour scheme is changed from the "Integral"
#     palette to the "Office Theme" palette (what the diff shows as the new
#     contents of the theme part driving the slide master).

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 16 -------------------------------------
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{4357684D-0C54-446E-AAA6-5CBF29CCA078}")
    }
}

# --- 2. Swap the theme colour scheme to the "Office Theme" palette ---------
function Set-ThemeRgb($colorItem, [string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $colorItem.RGB = $r + ($g * 256) + ($b * 65536)
}

$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

Set-ThemeRgb $colorScheme.Item(1)  "000000"   # dk1
Set-ThemeRgb $colorScheme.Item(2)  "FFFFFF"   # lt1
Set-ThemeRgb $colorScheme.Item(3)  "44546A"   # dk2
Set-ThemeRgb $colorScheme.Item(4)  "E7E6E6"   # lt2
Set-ThemeRgb $colorScheme.Item(5)  "5B9BD5"   # accent1
Set-ThemeRgb $colorScheme.Item(6)  "ED7D31"   # accent2
Set-ThemeRgb $colorScheme.Item(7)  "A5A5A5"   # accent3
Set-ThemeRgb $colorScheme.Item(8)  "FFC000"   # accent4
Set-ThemeRgb $colorScheme.Item(9)  "4472C4"   # accent5
Set-ThemeRgb $colorScheme.Item(10) "70AD47"   # accent6
Set-ThemeRgb $colorScheme.Item(11) "0563C1"   # hlink
Set-ThemeRgb $colorScheme.Item(12) "954F72"   # folHlink
